$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 140; this shifts rows 140:243 down to 141:244
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new data record.
$ws.Range("A140").Value = 9
$ws.Range("B140").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C140").Value = "Metropolitana"
$ws.Range("D140").Value = 44651
$ws.Range("E140").Value = 13
$ws.Range("F140").Value = 300000001
$ws.Range("G140").Value = "Rabanito"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 6100
$ws.Range("K140").Value = 3000
$ws.Range("L140").Value = 3000
$ws.Range("M140").Value = 3000
$ws.Range("N140").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O140").Value = "Provincia de Chacabuco"
$ws.Range("P140").Value = 30
$ws.Range("Q140").Value = 100
$ws.Range("R140").Value = "Hortaliza"
